# Remove embedded line breaks from the header row (row 3) labels on Sheet1,
# replacing the two-line "word<newline>unit" style headers with single-line
# text (using a space, or "(unit)" where the original used a bare "%" style).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "Sample No."
$ws.Range("D3").Value = "Wn %"
$ws.Range("E3").Value = "y (t/cu.m.)"
$ws.Range("F3").Value = "#4 (%)"
$ws.Range("G3").Value = "#200 (%)"
$ws.Range("H3").Value = "LL (%)"
$ws.Range("I3").Value = "PL (%)"
$ws.Range("J3").Value = "PI (%)"
$ws.Range("K3").Value = "Soil Class"
$ws.Range("L3").Value = "Su (t/sq.m.)"
$ws.Range("M3").Value = "Su from Ncor"

$ws.Range("S7").Select()
